$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.342.66'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '1.932.60'
$ws.Range("E3").Value = '  -2.27%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.607'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.30%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.47'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.90%  '
$ws.Range("E9").Value = '  -4.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0837'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.19%  '
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("D12").Value = '2.215.33'
$ws.Range("E12").Value = '  -2.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.798'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -7.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '13.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -11.58%  '
$ws.Range("E16").Value = '  -6.67%  '
$ws.Range("D17").Value = '1.943.70'
$ws.Range("E17").Value = '  -1.56%  '
$ws.Range("D18").Value = '36.254.45'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0859'
$ws.Range("E19").Value = '  -1.68%  '
$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '226.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.46%  '
$ws.Range("E22").Value = '  -7.27%  '
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("E24").Value = '  -11.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.18'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.130'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.10'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.72%  '
$ws.Range("E30").Value = '  -2.66%  '
$ws.Range("E31").Value = '  -6.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.52'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0622'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.06%  '
$ws.Range("E34").Value = '  -5.94%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.03'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("E38").Value = '  -6.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0961'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("E41").Value = '  -1.00%  '
$ws.Range("E42").Value = '  -3.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.19%  '
$ws.Range("D45").Value = '1.326.67'
$ws.Range("E45").Value = '  -2.99%  '
$ws.Range("E46").Value = '  -7.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.58%  '
$ws.Range("E49").Value = '  -0.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.95%  '
$ws.Range("D51").Value = '2.108.80'
$ws.Range("E51").Value = '  -2.20%  '
